# "Working on social vertex."
#
# Slide 3 is the bare socioeconomics/environmental-sociology/environmental-
# economics diagram (ovals + connectors + the three vertex labels). The
# author duplicated it to start building the next step of the diagram, so
# the duplicate (inserted right after slide 3) keeps the vertex labels while
# the original slide 3 has them stripped back off (hence "working on the
# social vertex" -- it's mid-edit on the new slide).

$p = $ppt.ActivePresentation

# Slide 3: the diagram-only slide (Oval 1 / Oval 9 / Oval 10 + connectors +
# the three vertex-label textboxes).
$diagramSlide = $p.Slides.Item(3)

# Duplicate it -- PowerPoint inserts the copy immediately after slide 3,
# pushing the old slide 4 down to slide 5. This mirrors the new
# p:sldId 260 appearing before the existing p:sldId 259 in sldIdLst.
$newSlide = $diagramSlide.Duplicate()

# Strip the three vertex-label textboxes back off the original slide 3,
# leaving just the ovals and connectors there.
$diagramSlide.Shapes.Item("TextBox 23").Delete()
$diagramSlide.Shapes.Item("TextBox 24").Delete()
$diagramSlide.Shapes.Item("TextBox 25").Delete()
